$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @{
    8  = -13.12309999999999
    10 = -13.5661
    12 = -10.4444
    18 = -11.5829
    37 = -12.8969
    55 = -13.4701
    68 = -11.3278
    77 = -13.64029999999999
    78 = -13.2715
    81 = -13.94289999999999
    82 = -11.9732
}

foreach ($row in $updates.Keys) {
    $ws.Range("C$row").Value = $updates[$row]
}

$wb.Save()
